$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell D3 value from 736 to 766
$ws.Range("D3").Value = 766

# Update the active selection to G8
$ws.Range("G8").Select()
